$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the formula in G8 to reference B8 instead of E8
$ws.Range("G8").Formula = "=86400*B8*0.3048^3"

# Add the new manually-calibrated value in H8
$ws.Range("H8").Value = 109028.77340000001

# Update the selected cell to reflect where the editor ended up
$ws.Range("G9").Select()
